# Update "想去人数" (F column) counts for a handful of events across
# the 展览 (Exhibition), 演出 (Performance), and 全部类型 (All types) sheets.

$wb = $excel.ActiveWorkbook

# 展览 sheet (sheetId=1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 913
$ws1.Range("F7").Value = 2643
$ws1.Range("F9").Value = 2579
$ws1.Range("F13").Value = 1677
$ws1.Range("F16").Value = 156
$ws1.Range("F25").Value = 578
$ws1.Range("F26").Value = 710
$ws1.Range("F31").Value = 1157
$ws1.Range("F32").Value = 189
$ws1.Range("F33").Value = 25
$ws1.Range("F34").Value = 1241
$ws1.Range("F36").Value = 299
$ws1.Range("F42").Value = 25

# 演出 sheet (sheetId=2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 74

# 全部类型 sheet (sheetId=4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 913
$ws4.Range("F6").Value = 2643
$ws4.Range("F7").Value = 2579
$ws4.Range("F8").Value = 1677
$ws4.Range("F13").Value = 156
$ws4.Range("F20").Value = 578
$ws4.Range("F21").Value = 710
$ws4.Range("F23").Value = 74
$ws4.Range("F29").Value = 1157
$ws4.Range("F30").Value = 189
$ws4.Range("F34").Value = 299
$ws4.Range("F42").Value = 25
